# Updated cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-50.
# NumberFormat is forced to "@" (Text) before assigning any Price value that would
# otherwise be auto-parsed as a number (e.g. "228.36"), so the cell keeps storing a
# literal text string just like the source data (prices with two dots, e.g.
# "34.613.54", are already unambiguous text and don't need the coercion guard).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.613.54"
$ws.Range("E2").Value = "  +1.15%  "

$ws.Range("D3").Value = "1.818.50"
$ws.Range("E3").Value = "  +1.79%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.36"
$ws.Range("E5").Value = "  +1.10%  "

$ws.Range("E6").Value = "  +1.05%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.79"
$ws.Range("E8").Value = "  +8.03%  "

$ws.Range("E9").Value = "  +2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0696"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("E11").Value = "  +0.38%  "

$ws.Range("D12").Value = "2.078.90"
$ws.Range("E12").Value = "  +1.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.38"
$ws.Range("E13").Value = "  +3.56%  "

$ws.Range("D14").Value = "1.829.82"
$ws.Range("E14").Value = "  +2.33%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.646"
$ws.Range("E15").Value = "  +3.28%  "

$ws.Range("D16").Value = "34.614.39"
$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("E17").Value = "  +3.76%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.18"
$ws.Range("E18").Value = "  +1.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.44"
$ws.Range("E19").Value = "  +0.62%  "

$ws.Range("E20").Value = "  +0.39%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.56"
$ws.Range("E21").Value = "  +5.93%  "

$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.22"
$ws.Range("E23").Value = "  +1.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "172.75"
$ws.Range("E24").Value = "  +6.91%  "

$ws.Range("E25").Value = "  +1.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.46"
$ws.Range("E26").Value = "  +4.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.77"
$ws.Range("E27").Value = "  +2.71%  "

$ws.Range("E28").Value = "  +1.22%  "

$ws.Range("E29").Value = "  -0.19%  "

$ws.Range("E30").Value = "  +7.73%  "

$ws.Range("E31").Value = "  +2.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").Value = "  +2.79%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.25"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("E34").Value = "  +2.94%  "

$ws.Range("D35").Value = "1.420.06"
$ws.Range("E35").Value = "  -1.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.59"
$ws.Range("E36").Value = "  -0.56%  "

$ws.Range("E37").Value = "  +2.30%  "

$ws.Range("E38").Value = "  +1.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0192"
$ws.Range("E39").Value = "  +1.12%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "85.98"
$ws.Range("E40").Value = "  +5.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.85"
$ws.Range("E41").Value = "  +4.38%  "

$ws.Range("E42").Value = "  +3.80%  "

$ws.Range("E43").Value = "  +0.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.81"
$ws.Range("E44").Value = "  -1.93%  "

$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("E46").Value = "  +2.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.11"
$ws.Range("E47").Value = "  +0.77%  "

$ws.Range("D48").Value = "1.981.15"
$ws.Range("E48").Value = "  +2.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.76"
$ws.Range("E49").Value = "  +0.43%  "

$ws.Range("E50").Value = "  +1.44%  "
